$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "visibility"
$ws.Range("E2").Value = $true
